$d = $word.ActiveDocument

$replacements = @(
    @("41×68=2788", "57×96=5472"),
    @("24×46=1104", "11×39=429"),
    @("32×56=1792", "59×74=4366"),
    @("63×65=4095", "31×90=2790"),
    @("41×74=3034", "34×94=3196"),
    @("31×15=465", "35×94=3290"),
    @("45×70=3150", "81×54=4374"),
    @("87×83=7221", "19×47=893"),
    @("50×22=1100", "73×18=1314"),
    @("76×89=6764", "68×42=2856"),
    @("47×30=1410", "32×50=1600"),
    @("22×32=704", "12×31=372"),
    @("88×21=1848", "91×74=6734"),
    @("64×26=1664", "57×99=5643"),
    @("19×21=399", "45×46=2070"),
    @("59×56=3304", "87×62=5394"),
    @("97×32=3104", "39×74=2886"),
    @("51×93=4743", "43×89=3827"),
    @("18×16=288", "83×96=7968"),
    @("62×91=5642", "68×56=3808"),
    @("67×67=4489", "14×41=574"),
    @("87×56=4872", "58×19=1102"),
    @("67×81=5427", "41×30=1230"),
    @("15×49=735", "43×15=645"),
    @("73×24=1752", "16×51=816")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
